# Update the header row (A1:D1) of the "Tabelle1" worksheet to the new
# lower-case / snake_case column names.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A1").Value = "jahr"
$ws.Range("B1").Value = "jahresschlusskurs"
$ws.Range("C1").Value = "ergebnis_pro_aktie"
$ws.Range("D1").Value = "dividendenvorschlag"

# Normalize C3's number format so it shares the same style as the other
# text-valued currency cells (C12, C13) instead of its own separate style.
$ws.Range("C3").NumberFormat = $ws.Range("C12").NumberFormat

# Leave the same cell selected as in the saved workbook.
$ws.Range("E11").Select()
